# Fills in the newly-known lesson dates for the "27.09 / 30.09" week block
# of the ОПІ-1 schedule table (table 1 on the page).
#
# Row with "30.09" (ПР04 / "Розроблення плану тестування ...") gains three
# dates in its previously-empty date columns: 05.10, 08.10, 12.10 - and the
# middle one (08.10) also gets the light-blue "accent5" cell shading that
# the other date-columns in the sheet use.
#
# The next two rows (Л04 / "Інфраструктура..." and ПР05 / "Визначення
# інфраструктури...") each get a single date filled into their first
# (currently empty) column: 04.10 and 07.10 respectively.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Run/paragraph-mark formatting shared by every cell in this table: Times
# New Roman, 14pt (half-point sz/szCs = 28).
$rPrXml = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'

function Set-CellDate($cell, [string]$text) {
    $openXml = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p><w:pPr>' + $rPrXml + '</w:pPr>' +
        '<w:r>' + $rPrXml + '<w:t>' + $text + '</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $cell.Range.InsertXML($openXml)
}

# --- Row with "30.09" in column 1: add 05.10 / 08.10 / 12.10 -------------
$row30 = 11
Set-CellDate $t.Cell($row30, 2) "05.10"
Set-CellDate $t.Cell($row30, 3) "08.10"
Set-CellDate $t.Cell($row30, 4) "12.10"

# Column 3 of that row also picks up the accent5 (light blue) cell shading.
$t.Cell($row30, 3).Shading.BackgroundPatternColor = 15986394  # RGB(218,238,243) = DAEEF3

# --- Row with "Л04": add 04.10 in column 1 --------------------------------
Set-CellDate $t.Cell(12, 1) "04.10"

# --- Row with "ПР05": add 07.10 in column 1 -------------------------------
Set-CellDate $t.Cell(13, 1) "07.10"
